$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.914.86"
$ws.Range("E2").Value = "  +6.96%  "

$ws.Range("D3").Value = "2.673.79"
$ws.Range("E3").Value = "  +10.18%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "511.34"
$ws.Range("E5").Value = "  +4.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.10"
$ws.Range("E6").Value = "  +2.77%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.606"
$ws.Range("E8").Value = "  +0.24%  "

$ws.Range("D9").Value = "2.670.88"
$ws.Range("E9").Value = "  +10.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.35"
$ws.Range("E10").Value = "  +7.66%  "

$ws.Range("E11").Value = "  +5.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.348"
$ws.Range("E12").Value = "  +4.03%  "

$ws.Range("E13").Value = "  +1.08%  "

$ws.Range("D14").Value = "3.135.71"
$ws.Range("E14").Value = "  +10.06%  "

$ws.Range("D15").Value = "60.969.14"
$ws.Range("E15").Value = "  +6.82%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.79"
$ws.Range("E16").Value = "  +4.95%  "

$ws.Range("E17").Value = "  +4.97%  "

$ws.Range("D18").Value = "2.669.82"
$ws.Range("E18").Value = "  +10.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.81"
$ws.Range("E19").Value = "  +1.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "347.59"
$ws.Range("E20").Value = "  +7.16%  "

$ws.Range("E21").Value = "  +4.97%  "

$ws.Range("E22").Value = "  +3.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.28"
$ws.Range("E24").Value = "  +3.56%  "

$ws.Range("E25").Value = "  +3.46%  "

$ws.Range("D26").Value = "2.771.23"
$ws.Range("E26").Value = "  +9.32%  "

$ws.Range("E27").Value = "  +3.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  -0.32%  "

$ws.Range("D29").Value = "0.0₃0860"
$ws.Range("E29").Value = "  +10.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.54"
$ws.Range("E30").Value = "  +3.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.48"
$ws.Range("E32").Value = "  +5.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.50"
$ws.Range("E33").Value = "  +5.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.71"
$ws.Range("E35").Value = "  +6.61%  "

$ws.Range("E36").Value = "  +8.69%  "

$ws.Range("E37").Value = "  +5.27%  "

$ws.Range("E38").Value = "  +11.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "310.51"
$ws.Range("E39").Value = "  +15.39%  "

$ws.Range("E40").Value = "  +1.54%  "

$ws.Range("E41").Value = "  +6.48%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.836"
$ws.Range("E42").Value = "  +28.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "35.40"
$ws.Range("E43").Value = "  +3.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.644"
$ws.Range("E44").Value = "  +8.63%  "

$ws.Range("E45").Value = "  +8.56%  "

$ws.Range("E46").Value = "  -0.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  -0.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.94"
$ws.Range("E48").Value = "  +14.76%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.85"
$ws.Range("E49").Value = "  +5.67%  "

$ws.Range("E50").Value = "  +3.64%  "

$ws.Range("D51").Value = "2.046.97"
$ws.Range("E51").Value = "  +9.57%  "
